# Re-colour the presentation's theme so the colour scheme that is currently
# "Integral / Red Violet" (stored in ppt/theme/theme1.xml, used by the one
# and only Slide Master) becomes the standard "Office Theme" colour scheme
# -- matching the authoring change that swapped the Office Theme content
# into ppt/theme/theme1.xml.
#
# PowerPoint's Theme.ThemeColorScheme exposes the 12 theme colour slots in
# a fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. Each slot's
# .RGB is a standard VB RGB() value (0x00BBGGRR) and setting it rewrites the
# corresponding <a:srgbClr val="…"/> entry in the theme part backing the
# active Slide Master.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# index -> (name, target "Office" RGB hex) just for reference/readability
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $scheme.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}

Write-Host "Theme colours updated to Office Theme palette."
